# Update the AMD Ryzen 5 5600X entry: switch it from the old DNS-shop
# listing to the new Citilink listing (new price + new product URL),
# and leave the cursor/selection on the updated cell (F2), matching the
# "Update notifications and shops" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New price for row 2 (AMD Ryzen 5 5600X)
$ws.Range("C2").Value = 8990

# New product URL text shown in column F for row 2
$ws.Range("F2").Value = "https://www.citilink.ru/product/processor-amd-ryzen-5-5600x-am4-oem-100-000000065-1773829/"

# Leave the active selection on F2, as in the saved workbook
$ws.Range("F2").Select()
